$wb = $excel.ActiveWorkbook

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 2267.5
$ws.Cells.Item(43, 9).Value = 2000
$ws.Cells.Item(43, 10).Value = 2891.6667
$ws.Cells.Item(43, 11).Value = 2000
$ws.Cells.Item(43, 12).Value = 2891.6667
$ws.Cells.Item(43, 13).Value = -1931
$ws.Cells.Item(43, 14).Value = -3029.6667

# ALC row 59
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(59, 8).Value = 2000
$ws.Cells.Item(59, 10).Value = 2000
$ws.Cells.Item(59, 12).Value = 6000
$ws.Cells.Item(59, 14).Value = -7114

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 4359
$ws.Cells.Item(76, 9).Value = 4323.75
$ws.Cells.Item(76, 10).Value = 4500
$ws.Cells.Item(76, 11).Value = 4323.75
$ws.Cells.Item(76, 12).Value = 4500
$ws.Cells.Item(76, 13).Value = -4008.75
$ws.Cells.Item(76, 14).Value = -5130

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 8).Value = 4359
$ws.Cells.Item(79, 9).Value = 4323.75
$ws.Cells.Item(79, 10).Value = 4500
$ws.Cells.Item(79, 11).Value = 4323.75
$ws.Cells.Item(79, 12).Value = 4500
$ws.Cells.Item(79, 13).Value = -3231.75
$ws.Cells.Item(79, 14).Value = -6684

# ALC row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(118, 8).Value = 1084.6364
$ws.Cells.Item(118, 10).Value = 1748
$ws.Cells.Item(118, 12).Value = 5244
$ws.Cells.Item(118, 14).Value = -8558

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 10559.8
$ws.Cells.Item(125, 9).Value = 5000
$ws.Cells.Item(125, 11).Value = 45000
$ws.Cells.Item(125, 13).Value = -42540

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1450.8
$ws.Cells.Item(137, 9).Value = 1625.375
$ws.Cells.Item(137, 10).Value = 1251.2858
$ws.Cells.Item(137, 11).Value = 4876.125
$ws.Cells.Item(137, 12).Value = 3753.8574
$ws.Cells.Item(137, 13).Value = -2326.125
$ws.Cells.Item(137, 14).Value = -8853.857400000001

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2089.6667
$ws.Cells.Item(2, 9).Value = 2108.1
$ws.Cells.Item(2, 11).Value = 2108.1
$ws.Cells.Item(2, 13).Value = -1995.1

# ARM row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 939.3333
$ws.Cells.Item(4, 9).Value = 956.75
$ws.Cells.Item(4, 10).Value = 800
$ws.Cells.Item(4, 11).Value = 956.75
$ws.Cells.Item(4, 12).Value = 800
$ws.Cells.Item(4, 13).Value = -840.75
$ws.Cells.Item(4, 14).Value = -1032

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 61.166668
$ws.Cells.Item(5, 9).Value = 42.5
$ws.Cells.Item(5, 11).Value = 42.5
$ws.Cells.Item(5, 13).Value = 69.5

# ARM row 53
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 13).Value = $null

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1933.3334
$ws.Cells.Item(61, 9).Value = 1933.3334
$ws.Cells.Item(61, 11).Value = 1933.3334
$ws.Cells.Item(61, 13).Value = -1721.3334

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 3956
$ws.Cells.Item(74, 9).Value = 3956
$ws.Cells.Item(74, 11).Value = 3956
$ws.Cells.Item(74, 13).Value = -3082

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 3956
$ws.Cells.Item(77, 9).Value = 3956
$ws.Cells.Item(77, 11).Value = 19780
$ws.Cells.Item(77, 13).Value = -15412

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 2089.6667
$ws.Cells.Item(116, 9).Value = 2108.1
$ws.Cells.Item(116, 11).Value = 2108.1
$ws.Cells.Item(116, 13).Value = 185.9000000000001

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1933.3334
$ws.Cells.Item(136, 9).Value = 1933.3334
$ws.Cells.Item(136, 11).Value = 5800.0002
$ws.Cells.Item(136, 13).Value = -3250.0002

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2089.6667
$ws.Cells.Item(3, 9).Value = 2108.1
$ws.Cells.Item(3, 11).Value = 2108.1
$ws.Cells.Item(3, 13).Value = -1994.1

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 61.166668
$ws.Cells.Item(4, 9).Value = 42.5
$ws.Cells.Item(4, 11).Value = 42.5
$ws.Cells.Item(4, 13).Value = 72.5

# BSM row 68
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(68, 8).Value = 25000
$ws.Cells.Item(68, 9).Value = 25000
$ws.Cells.Item(68, 11).Value = 25000
$ws.Cells.Item(68, 13).Value = -24189

# BSM row 71
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(71, 8).Value = 25000
$ws.Cells.Item(71, 9).Value = 25000
$ws.Cells.Item(71, 11).Value = 75000
$ws.Cells.Item(71, 13).Value = -70944

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2273
$ws.Cells.Item(86, 9).Value = 1325.5625
$ws.Cells.Item(86, 10).Value = 7326
$ws.Cells.Item(86, 11).Value = 1325.5625
$ws.Cells.Item(86, 12).Value = 7326
$ws.Cells.Item(86, 13).Value = -202.5625
$ws.Cells.Item(86, 14).Value = -9572

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 2273
$ws.Cells.Item(89, 9).Value = 1325.5625
$ws.Cells.Item(89, 10).Value = 7326
$ws.Cells.Item(89, 11).Value = 6627.8125
$ws.Cells.Item(89, 12).Value = 36630
$ws.Cells.Item(89, 13).Value = -1011.8125
$ws.Cells.Item(89, 14).Value = -47862

# BSM row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(135, 8).Value = 47250
$ws.Cells.Item(135, 10).Value = 47250
$ws.Cells.Item(135, 12).Value = 47250
$ws.Cells.Item(135, 14).Value = -57390

# CRP row 2
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 324.66666
$ws.Cells.Item(2, 9).Value = 305.75
$ws.Cells.Item(2, 10).Value = 362.5
$ws.Cells.Item(2, 11).Value = 305.75
$ws.Cells.Item(2, 12).Value = 362.5
$ws.Cells.Item(2, 13).Value = -192.75
$ws.Cells.Item(2, 14).Value = -588.5

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 250.44444
$ws.Cells.Item(22, 9).Value = 230.4
$ws.Cells.Item(22, 10).Value = 275.5
$ws.Cells.Item(22, 11).Value = 230.4
$ws.Cells.Item(22, 12).Value = 275.5
$ws.Cells.Item(22, 13).Value = 119.6
$ws.Cells.Item(22, 14).Value = -975.5

# CRP row 35
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(35, 8).Value = 7436.3335
$ws.Cells.Item(35, 10).Value = 20000
$ws.Cells.Item(35, 12).Value = 20000
$ws.Cells.Item(35, 14).Value = -20588

# CRP row 97
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(97, 8).Value = 80000
$ws.Cells.Item(97, 10).Value = 80000
$ws.Cells.Item(97, 12).Value = 80000
$ws.Cells.Item(97, 14).Value = -81982

# CRP row 109
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(109, 8).Value = 46606.668
$ws.Cells.Item(109, 9).Value = 53259
$ws.Cells.Item(109, 10).Value = 43280.5
$ws.Cells.Item(109, 11).Value = 53259
$ws.Cells.Item(109, 12).Value = 43280.5
$ws.Cells.Item(109, 13).Value = -52219
$ws.Cells.Item(109, 14).Value = -45360.5

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 260.4762
$ws.Cells.Item(2, 9).Value = 52.875
$ws.Cells.Item(2, 10).Value = 924.8
$ws.Cells.Item(2, 11).Value = 317.25
$ws.Cells.Item(2, 12).Value = 5548.799999999999
$ws.Cells.Item(2, 13).Value = -204.25
$ws.Cells.Item(2, 14).Value = -5774.799999999999

# CUL row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 561.4
$ws.Cells.Item(7, 9).Value = 252.5
$ws.Cells.Item(7, 10).Value = 767.3333
$ws.Cells.Item(7, 11).Value = 757.5
$ws.Cells.Item(7, 12).Value = 2301.9999
$ws.Cells.Item(7, 13).Value = -645.5
$ws.Cells.Item(7, 14).Value = -2525.9999

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 1499.8334
$ws.Cells.Item(55, 9).Value = 1000
$ws.Cells.Item(55, 11).Value = 3000
$ws.Cells.Item(55, 13).Value = -2823

# CUL row 124
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 14).Value = $null

# CUL row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 13).Value = $null

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 524.7143
$ws.Cells.Item(2, 9).Value = 515
$ws.Cells.Item(2, 11).Value = 515
$ws.Cells.Item(2, 13).Value = -402

# GSM row 123
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 59999
$ws.Cells.Item(123, 10).Value = 59999
$ws.Cells.Item(123, 12).Value = 59999
$ws.Cells.Item(123, 14).Value = -64899

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1999.8
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 1999.8
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 1999.8
$ws.Cells.Item(22, 13).Value = $null
$ws.Cells.Item(22, 14).Value = -2589.8

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 1999.8
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 1999.8
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 1999.8
$ws.Cells.Item(27, 13).Value = $null
$ws.Cells.Item(27, 14).Value = -2213.8

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2277.652
$ws.Cells.Item(46, 9).Value = 2069.35
$ws.Cells.Item(46, 10).Value = 3666.3333
$ws.Cells.Item(46, 11).Value = 2069.35
$ws.Cells.Item(46, 12).Value = 3666.3333
$ws.Cells.Item(46, 13).Value = -1881.35
$ws.Cells.Item(46, 14).Value = -4042.3333

# LTW row 57
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(57, 8).Value = 8000
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 13).Value = $null

# LTW row 123
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(123, 8).Value = 78499.5
$ws.Cells.Item(123, 10).Value = 78499.5
$ws.Cells.Item(123, 12).Value = 78499.5
$ws.Cells.Item(123, 14).Value = -88299.5

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 16547.5
$ws.Cells.Item(132, 9).Value = 16547.5
$ws.Cells.Item(132, 11).Value = 49642.5
$ws.Cells.Item(132, 13).Value = -47112.5

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(133, 8).Value = 49999
$ws.Cells.Item(133, 9).Value = 49999
$ws.Cells.Item(133, 11).Value = 49999
$ws.Cells.Item(133, 13).Value = -47469

# WVR row 69
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 9723
$ws.Cells.Item(69, 10).Value = 9723
$ws.Cells.Item(69, 12).Value = 9723
$ws.Cells.Item(69, 14).Value = -11221

# WVR row 72
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(72, 8).Value = 9723
$ws.Cells.Item(72, 10).Value = 9723
$ws.Cells.Item(72, 12).Value = 29169
$ws.Cells.Item(72, 14).Value = -36657

# WVR row 109
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 67999.5
$ws.Cells.Item(109, 10).Value = 67999.5
$ws.Cells.Item(109, 12).Value = 67999.5
$ws.Cells.Item(109, 14).Value = -70773.5

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1922.3334
$ws.Cells.Item(136, 9).Value = 1900.125
$ws.Cells.Item(136, 11).Value = 5700.375
$ws.Cells.Item(136, 13).Value = -3150.375
